$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append after the existing last data row (187: 28-09-2021).
# Column A holds the date as plain text (matches the existing "dd-mm-yyyy"
# string series); columns B/C/D repeat the current TPM / facilidades values.
$newDates = @("29-09-2021", "30-09-2021", "01-10-2021", "04-10-2021", "05-10-2021")

$startRow = 188
for ($i = 0; $i -lt $newDates.Count; $i++) {
    $r = $startRow + $i
    $dateCell = $ws.Cells.Item($r, 1)

    # Some of these strings (e.g. "01-10-2021") are ambiguous and would be
    # silently reinterpreted as a date serial number by a plain
    # `.Value = "..."` assignment. Route the text through a formula result
    # (always text-typed) and then freeze it to a static value via
    # copy / paste-special-values, so the cell keeps its default (General)
    # style and ends up as plain shared-string text, exactly like the rest
    # of the column.
    $dateCell.Formula = '="' + $newDates[$i] + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)  # xlPasteValues

    $ws.Cells.Item($r, 2).Value = 1.5
    $ws.Cells.Item($r, 3).Value = 1.75
    $ws.Cells.Item($r, 4).Value = 1.25
}

$excel.CutCopyMode = 0
